# 10th - Bug Fixes & Tickers Update
# Insert a new "Jun_26" date column right after the firm/UN columns (before
# the old "Jun_17" column), fill it in, and append two new tracked firms.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank column before column B. This shifts the existing
#    Jun_17 / Jun_15 / Jun_13 / Jun_10 columns (B:E) right to (C:F).
$ws.Columns("B:B").Insert()

# 2. Header for the freshly inserted column: the newest snapshot date.
$ws.Range("B1").Value = "Jun_26"

# 3. Match the column's width to its neighbors.
$ws.Columns("B:B").ColumnWidth = 56.83203125

# 4. Fill the new column's data rows with "UN" (unchanged), same as every
#    other firm/date cell that saw no rating action on this date.
for ($row = 2; $row -le 27; $row++) {
    $ws.Cells.Item($row, 2).Value = "UN"
}

# 5. Append two newly tracked research firms at the bottom of the table.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
